$d = $word.ActiveDocument

# Insert the new paragraphs at the very end of the document body (after the
# last paragraph mark), so the existing "_GoBack" bookmark in the
# "Simply move the mouse..." paragraph is left untouched by this step.
$endPos = $d.Content.End
$insertRange = $d.Range($endPos, $endPos)

$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Since the computer's power is on, we will have to acquire the volatile data.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>These are data that are lost once the computer is shut down.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">The following are the kind of volatile data you will want to gather: system time, RAM, process information, network log, logged-on </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>users’</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> information, and cached data (command history, clipboard, print spool files).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Please ensure you gather the data in the order of most to least volatile. The most volatile data </w:t></w:r><w:r><w:t>are</w:t></w:r><w:r><w:t xml:space="preserve"> the ones that change consistently as time goes on, and the least being the ones that </w:t></w:r><w:r><w:t>rarely change</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@

$insertRange.InsertXML($fragment) | Out-Null

# Relocate the "_GoBack" bookmark: it used to sit right after "Simply move
# the mouse..." (the end of that paragraph); it now belongs at the very
# start of the final paragraph ("Please ensure you gather the data...").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
